# Auto-generated edit script: updates Kraken_Profits market-price snapshot values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 165.9
$ws.Range("I9").Value = 182.44444
$ws.Range("K9").Value = 182.44444
$ws.Range("M9").Value = -13.44443999999999

$ws.Range("H17").Value = 2098.8462
$ws.Range("J17").Value = 2080.4546
$ws.Range("L17").Value = 6241.3638
$ws.Range("N17").Value = -6577.3638

$ws.Range("H18").Value = 3199.8
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5568

$ws.Range("H51").Value = 12000
$ws.Range("J51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("N51").Value = -12968

$ws.Range("H131").Value = 781
$ws.Range("I131").Value = 781
$ws.Range("K131").Value = 2343
$ws.Range("M131").Value = 2697

$ws.Range("H137").Value = 1817.8
$ws.Range("I137").Value = 1779.8334
$ws.Range("K137").Value = 5339.5002
$ws.Range("M137").Value = -2789.5002

$ws.Range("H138").Value = 2829.16
$ws.Range("I138").Value = 749.44446
$ws.Range("K138").Value = 2248.33338
$ws.Range("M138").Value = 2891.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5221.3
$ws.Range("I32").Value = 5221.3
$ws.Range("K32").Value = 5221.3
$ws.Range("M32").Value = -4934.3

$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4671

$ws.Range("H36").Value = 3675.3333
$ws.Range("I36").Value = 3675.3333
$ws.Range("K36").Value = 3675.3333
$ws.Range("M36").Value = -3329.3333

$ws.Range("H61").Value = 4503.75
$ws.Range("I61").Value = 4660.75
$ws.Range("J61").Value = 4346.75
$ws.Range("K61").Value = 4660.75
$ws.Range("L61").Value = 4346.75
$ws.Range("M61").Value = -4448.75
$ws.Range("N61").Value = -4770.75

$ws.Range("H132").Value = 4200
$ws.Range("I132").Value = 4200
$ws.Range("K132").Value = 12600
$ws.Range("M132").Value = -10070

$ws.Range("H136").Value = 4503.75
$ws.Range("I136").Value = 4660.75
$ws.Range("J136").Value = 4346.75
$ws.Range("K136").Value = 13982.25
$ws.Range("L136").Value = 13040.25
$ws.Range("M136").Value = -11432.25
$ws.Range("N136").Value = -18140.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12337
$ws.Range("I82").Value = 12337
$ws.Range("K82").Value = 12337
$ws.Range("M82").Value = -11954

$ws.Range("H85").Value = 12337
$ws.Range("I85").Value = 12337
$ws.Range("K85").Value = 12337
$ws.Range("M85").Value = -11011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5037.3076
$ws.Range("I31").Value = 4311.5
$ws.Range("K31").Value = 4311.5
$ws.Range("M31").Value = -4016.5

$ws.Range("H34").Value = 5037.3076
$ws.Range("I34").Value = 4311.5
$ws.Range("K34").Value = 4311.5
$ws.Range("M34").Value = -4109.5

$ws.Range("H62").Value = 4174.75
$ws.Range("I62").Value = 3899.6667
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3899.6667
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3275.6667
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 4174.75
$ws.Range("I65").Value = 3899.6667
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 19498.3335
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16378.3335
$ws.Range("N65").Value = -31240

$ws.Range("H86").Value = 4749
$ws.Range("I86").Value = 4749
$ws.Range("K86").Value = 4749
$ws.Range("M86").Value = -3626

$ws.Range("H89").Value = 4749
$ws.Range("I89").Value = 4749
$ws.Range("K89").Value = 23745
$ws.Range("M89").Value = -18129

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 499.75
$ws.Range("I11").Value = 500.5
$ws.Range("J11").Value = 499
$ws.Range("K11").Value = 1501.5
$ws.Range("L11").Value = 1497
$ws.Range("M11").Value = -1361.5
$ws.Range("N11").Value = -1777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12375

$ws.Range("H97").Value = 487.7143
$ws.Range("I97").Value = 487.7143
$ws.Range("K97").Value = 487.7143
$ws.Range("M97").Value = 8.28570000000002

$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2699
$ws.Range("I7").Value = 2699
$ws.Range("K7").Value = 2699
$ws.Range("M7").Value = -2587

$ws.Range("H68").Value = 2999.8333
$ws.Range("I68").Value = 2999.8
$ws.Range("K68").Value = 2999.8
$ws.Range("M68").Value = -2250.8

$ws.Range("H71").Value = 2999.8333
$ws.Range("I71").Value = 2999.8
$ws.Range("K71").Value = 14999
$ws.Range("M71").Value = -11255

$ws.Range("H126").Value = 2699
$ws.Range("I126").Value = 2699
$ws.Range("K126").Value = 8097
$ws.Range("M126").Value = -5627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 40985.2
$ws.Range("I4").Value = 50544
$ws.Range("J4").Value = 2750
$ws.Range("K4").Value = 50544
$ws.Range("L4").Value = 2750
$ws.Range("M4").Value = -50431
$ws.Range("N4").Value = -2976

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H81").Value = 31800.4
$ws.Range("I81").Value = 31800.4
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 63600.8
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -62539.8
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 31800.4
$ws.Range("I84").Value = 31800.4
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 318004
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -312700
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 712.3333
$ws.Range("I113").Value = 769.25
$ws.Range("K113").Value = 2307.75
$ws.Range("M113").Value = -137.75

$ws.Range("H136").Value = 2502
$ws.Range("I136").Value = 2502
$ws.Range("K136").Value = 7506
$ws.Range("M136").Value = -4956
